# Fill in the bill form with teacher/term details and the amount in words.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = "নিয়মিত পরীক্ষা ২০২২"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"
$ws.Range("B5").Value = "সিএসই"
$ws.Range("A3").Value = "নাম: Ms. Dola Das "
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :সিএসই"
$ws.Range("A32").Value = "কথায়:চার লক্ষ একচল্লিশ হাজার বত্রিশ টাকা মাত্র।"

# Column A is widened to fit the longer "নাম: ..." text.
$ws.Columns.Item(1).ColumnWidth = 14.33203125

# Row 36 grows to accommodate wrapped signature text.
$ws.Rows.Item(36).RowHeight = 68.4

# Scroll/selection state to match the end of the editing session.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I32").Select()
